$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '66.612.03'
$ws.Range('E2').Value = '  +5.17%  '
$ws.Range('D3').Value = '3.729.40'
$ws.Range('E3').Value = '  +7.70%  '
$ws.Range('E4').Value = '  +0.20%  '
Set-TextValue $ws.Range('D5') '422.87'
$ws.Range('E5').Value = '  +1.80%  '
Set-TextValue $ws.Range('D6') '131.68'
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('D7').Value = '3.723.32'
$ws.Range('E7').Value = '  +7.84%  '
Set-TextValue $ws.Range('D8') '0.640'
$ws.Range('E8').Value = '  +2.11%  '
Set-TextValue $ws.Range('D9') '0.999'
$ws.Range('E9').Value = '  -0.09%  '
Set-TextValue $ws.Range('D10') '0.762'
$ws.Range('E10').Value = '  -1.01%  '
Set-TextValue $ws.Range('D11') '0.180'
$ws.Range('E11').Value = '  +14.35%  '
Set-TextValue $ws.Range('D12') '0.0000375'
$ws.Range('E12').Value = '  +61.10%  '
Set-TextValue $ws.Range('D13') '42.58'
$ws.Range('E13').Value = '  +0.98%  '
Set-TextValue $ws.Range('D14') '10.25'
$ws.Range('E14').Value = '  +4.50%  '
$ws.Range('D15').Value = '4.297.48'
$ws.Range('E15').Value = '  +6.87%  '
$ws.Range('E16').Value = '  -0.16%  '
Set-TextValue $ws.Range('D17') '20.95'
$ws.Range('E17').Value = '  +3.19%  '
$ws.Range('D18').Value = '3.704.82'
$ws.Range('E18').Value = '  +7.16%  '
Set-TextValue $ws.Range('D19') '12.99'
$ws.Range('E19').Value = '  +4.37%  '
$ws.Range('E20').Value = '  +2.92%  '
$ws.Range('D21').Value = '66.702.18'
$ws.Range('E21').Value = '  +5.43%  '
Set-TextValue $ws.Range('D22') '445.94'
$ws.Range('E22').Value = '  -2.64%  '
Set-TextValue $ws.Range('D23') '15.67'
$ws.Range('E23').Value = '  +16.42%  '
Set-TextValue $ws.Range('D24') '89.81'
$ws.Range('E24').Value = '  -0.61%  '
Set-TextValue $ws.Range('D25') '3.18'
$ws.Range('E25').Value = '  -3.55%  '
Set-TextValue $ws.Range('D26') '38.04'
$ws.Range('E26').Value = '  +13.18%  '
Set-TextValue $ws.Range('D27') '10.19'
$ws.Range('E27').Value = '  +0.58%  '
Set-TextValue $ws.Range('D28') '3.29'
$ws.Range('E28').Value = '  -0.22%  '
Set-TextValue $ws.Range('D29') '4.99'
$ws.Range('E29').Value = '  +4.93%  '
Set-TextValue $ws.Range('D30') '12.56'
$ws.Range('E30').Value = '  +1.38%  '
$ws.Range('E31').Value = '  +3.82%  '
Set-TextValue $ws.Range('D32') '0.120'
$ws.Range('E32').Value = '  +7.06%  '
Set-TextValue $ws.Range('D33') '7.27'
$ws.Range('E33').Value = '  -3.18%  '
Set-TextValue $ws.Range('D34') '42.10'
$ws.Range('E34').Value = '  +5.30%  '
$ws.Range('E35').Value = '  -1.77%  '
$ws.Range('E36').Value = '  -0.02%  '
Set-TextValue $ws.Range('D37') '56.50'
$ws.Range('E37').Value = '  -2.08%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').Value = '0.0₃0730'
$ws.Range('E39').Value = '  +13.86%  '
Set-TextValue $ws.Range('D40') '3.03'
$ws.Range('E40').Value = '  +30.35%  '
Set-TextValue $ws.Range('D41') '0.146'
$ws.Range('E41').Value = '  +5.61%  '
Set-TextValue $ws.Range('D42') '28.53'
$ws.Range('E42').Value = '  +30.44%  '
$ws.Range('B43').Value = 'LidoDAOToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D43') '3.49'
$ws.Range('E43').Value = '  +4.44%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D44') '0.998'
$ws.Range('E44').Value = '  -0.08%  '
Set-TextValue $ws.Range('D45') '2.16'
$ws.Range('E45').Value = '  +8.36%  '
Set-TextValue $ws.Range('D46') '145.60'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('E47').Value = '  -5.71%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D48') '4.37'
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D49') '2.65'
$ws.Range('E49').Value = '  -5.88%  '
$ws.Range('E50').Value = '  -4.09%  '
$ws.Range('E51').Value = '  +14.67%  '
